$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.281000000000003
$ws.Range("B14").Value = 5.504300000000002
$ws.Range("B16").Value = 6.284200000000001
$ws.Range("B21").Value = 9.529500000000002
$ws.Range("B23").Value = 9.112600000000002
$ws.Range("B25").Value = 5.994799999999998
$ws.Range("B26").Value = 5.48220000000001
$ws.Range("B29").Value = 4.982200000000003
$ws.Range("B40").Value = 8.843700000000002
$ws.Range("B53").Value = 5.200799999999998
$ws.Range("B57").Value = 4.893999999999998
$ws.Range("B59").Value = 4.840299999999998
$ws.Range("B65").Value = 6.058200000000003
$ws.Range("B69").Value = 5.425899999999996
$ws.Range("B79").Value = 9.196400000000002
$ws.Range("B83").Value = 5.8766
$ws.Range("B91").Value = 5.134900000000001
$ws.Range("B93").Value = 5.5629
$ws.Range("B100").Value = 5.650999999999999
